# Apply crypto price/volume updates to match target snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.286.66'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.838.40'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.70%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.09%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '238.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.68%  '

$ws.Range("E6").Value = '  -1.07%  '

$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07371'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.40%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2886'
$ws.Range("D9").Style = "Normal"

$ws.Range("E10").Value = '  -1.54%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07721'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.39%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.837.81'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.82%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.945'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.81%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.00001050'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.95%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6623'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.09%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '81.30'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.82%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.261'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.26%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.265.71'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.15%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '235.21'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +2.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.22'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.55%  '

$ws.Range("E21").Value = '  +0.15%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.268'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.54%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.002'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.15%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '157.23'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.49%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.409'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.36%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1332'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.63%  '

$ws.Range("E27").Value = '  -1.93%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.07076'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.76%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.474'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.26%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.480'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.56%  '

$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.022'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -2.19%  '

$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.014'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.43%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.148'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.26%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.786'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6921'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.35%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.586'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.71%  '

$ws.Range("E37").Value = '  -2.47%  '

$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.233.61'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.07%  '

$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.780'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.18%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.759'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9461'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.53%  '

$ws.Range("E42").Value = '  -0.01%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.989.91'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.67%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.18'
$ws.Range("D44").Style = "Normal"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '65.13'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.72%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000119'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.30%  '

$ws.Range("E47").Value = '  -2.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.680'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.30%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '8.899'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.80%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.1126'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -3.16%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.3873'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.06%  '
